$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("protocol_info")

# New set of protocol-code headers for columns B..Z (row 1), replacing the
# old B..N set and extending the table out to column Z.
$headerVals = @("0046","0052","0048","0051","0037","0067","0054","0053","0068","0044","0050","A251","0059","0055","0062","0047","0079I","0049","0064","0071","0073","0069","0069I","0075","0075S")

# CTN-NODES / CTN-SITES row counts for the (now 25-wide) protocol columns.
$row2Vals = @(6,4,5,8,3,4,2,6,5,2,2,1,3,2,2,2,1,10,7,1,1,4,4,1,1)
$row3Vals = @(8,4,6,9,3,8,3,8,9,3,3,1,5,2,6,2,1,11,11,1,1,4,4,3,3)

$firstCol = 2   # column B
$lastCol  = 26  # column Z

# --- Row 1 header text -------------------------------------------------
# Write each header as a formula string first so the engine treats it as
# text even when it looks purely numeric (e.g. "0046"), then convert the
# whole row to static values and finally restore the bold/bordered header
# style that already lives on column B's header cell.
for ($c = $firstCol; $c -le $lastCol; $c++) {
    $val = $headerVals[$c - $firstCol]
    $ws.Cells.Item(1, $c).Formula = '="' + $val + '"'
}
$headerRange = $ws.Range($ws.Cells.Item(1, $firstCol), $ws.Cells.Item(1, $lastCol))
$headerRange.Copy()
$headerRange.PasteSpecial(-4163)  # xlPasteValues
$ws.Cells.Item(1, 2).Copy()       # sample header style (bold, centered, bordered)
$headerRange.PasteSpecial(-4122)  # xlPasteFormats

# --- Row labels (A2 / A3) ----------------------------------------------
$ws.Cells.Item(2, 1).Formula = '="CTN-NODES"'
$ws.Cells.Item(3, 1).Formula = '="CTN-SITES"'
$labelRange = $ws.Range($ws.Cells.Item(2, 1), $ws.Cells.Item(3, 1))
$labelRange.Copy()
$labelRange.PasteSpecial(-4163)
$ws.Cells.Item(2, 1).Copy()
$labelRange.PasteSpecial(-4122)

# --- Numeric data rows 2 & 3 --------------------------------------------
for ($c = $firstCol; $c -le $lastCol; $c++) {
    $ws.Cells.Item(2, $c).Value = $row2Vals[$c - $firstCol]
    $ws.Cells.Item(3, $c).Value = $row3Vals[$c - $firstCol]
}

$excel.CutCopyMode = $false
